$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I8").Value = 274
$ws.Range("I10").Value = 962
$ws.Range("I11").Value = 540
$ws.Range("I20").Value = 960
$ws.Range("I22").Value = 1040
$ws.Range("I23").Value = 347
$ws.Range("I29").Value = 787
$ws.Range("I43").Value = 646
$ws.Range("I47").Value = 934
$ws.Range("I49").Value = 719
$ws.Range("I51").Value = 486
$ws.Range("I55").Value = 348
$ws.Range("I57").Value = 1006
$ws.Range("I65").Value = 391
$ws.Range("I68").Value = 660
$ws.Range("I76").Value = 582
$ws.Range("I77").Value = 985
$ws.Range("I82").Value = 899
$ws.Range("I88").Value = 1007
$ws.Range("I93").Value = 1039
$ws.Range("I109").Value = 876
$ws.Range("I110").Value = 702
$ws.Range("I114").Value = 378
$ws.Range("I120").Value = 443
$ws.Range("I122").Value = 442
$ws.Range("I123").Value = 817
$ws.Range("I125").Value = 437
$ws.Range("I131").Value = 701
$ws.Range("I132").Value = 902
$ws.Range("I139").Value = 889
$ws.Range("I143").Value = 877
$ws.Range("I145").Value = 740
$ws.Range("I146").Value = 803
$ws.Range("I149").Value = 1029
$ws.Range("I151").Value = 949
$ws.Range("I157").Value = 993
$ws.Range("I162").Value = 711
$ws.Range("I164").Value = 758
$ws.Range("I168").Value = 604
$ws.Range("I171").Value = 727
$ws.Range("I174").Value = 527
$ws.Range("I175").Value = 779
$ws.Range("I180").Value = 976
$ws.Range("I185").Value = 872
$ws.Range("I190").Value = 539
$ws.Range("I191").Value = 819
$ws.Range("I192").Value = 997
$ws.Range("I213").Value = 990
$ws.Range("I214").Value = 802
$ws.Range("I215").Value = 682
$ws.Range("I216").Value = 746
$ws.Range("I221").Value = 907
$ws.Range("I222").Value = 1026
$ws.Range("I236").Value = 887
$ws.Range("I237").Value = 744
$ws.Range("I238").Value = 536
$ws.Range("I247").Value = 273
$ws.Range("I249").Value = 369
$ws.Range("I256").Value = 859
$ws.Range("I257").Value = 583
$ws.Range("I269").Value = 1051
$ws.Range("I273").Value = 823
$ws.Range("I276").Value = 1033
$ws.Range("I277").Value = 968
$ws.Range("I285").Value = 795
$ws.Range("I292").Value = 485
$ws.Range("I299").Value = 905
$ws.Range("I300").Value = 510
$ws.Range("I302").Value = 647
$ws.Range("I308").Value = 567
$ws.Range("I310").Value = 810
$ws.Range("I320").Value = 942
$ws.Range("I321").Value = 669
$ws.Range("I334").Value = 747
$ws.Range("I342").Value = 868
$ws.Range("I352").Value = 980
$ws.Range("I355").Value = 836
$ws.Range("I360").Value = 829
$ws.Range("I361").Value = 640
$ws.Range("I366").Value = 816
$ws.Range("I379").Value = 1015
$ws.Range("I386").Value = 718
$ws.Range("I395").Value = 594
$ws.Range("I396").Value = 900
$ws.Range("I401").Value = 710
$ws.Range("I403").Value = 866
$ws.Range("I422").Value = 603
$ws.Range("I425").Value = 505
$ws.Range("I435").Value = 834
$ws.Range("I436").Value = 996
$ws.Range("I443").Value = 837
$ws.Range("I444").Value = 792
$ws.Range("I447").Value = 906
$ws.Range("I448").Value = 537
$ws.Range("I456").Value = 820
$ws.Range("I464").Value = 504
$ws.Range("I467").Value = 613
$ws.Range("I469").Value = 835
$ws.Range("I472").Value = 392
$ws.Range("I475").Value = 757
$ws.Range("I476").Value = 528
$ws.Range("I481").Value = 487
$ws.Range("I483").Value = 659
$ws.Range("I485").Value = 601
$ws.Range("I488").Value = 316
$ws.Range("I492").Value = 733
$ws.Range("I495").Value = 788
$ws.Range("I497").Value = 568
$ws.Range("I501").Value = 918
$ws.Range("I506").Value = 975
$ws.Range("I508").Value = 861
$ws.Range("I511").Value = 964
$ws.Range("I514").Value = 748
$ws.Range("I529").Value = 914
$ws.Range("I536").Value = 670
$ws.Range("I538").Value = 941
$ws.Range("I546").Value = 703
$ws.Range("I547").Value = 979
$ws.Range("I555").Value = 954
$ws.Range("I557").Value = 488
$ws.Range("I561").Value = 974
$ws.Range("I570").Value = 958
$ws.Range("I571").Value = 1056
$ws.Range("I572").Value = 370
$ws.Range("I575").Value = 1063
$ws.Range("I578").Value = 563
$ws.Range("I583").Value = 967
$ws.Range("I587").Value = 860
$ws.Range("I597").Value = 924
$ws.Range("I605").Value = 675
$ws.Range("I619").Value = 674
$ws.Range("I625").Value = 982
$ws.Range("I628").Value = 963
$ws.Range("I631").Value = 973
$ws.Range("I634").Value = 981
$ws.Range("I640").Value = 436
$ws.Range("I659").Value = 1048
$ws.Range("I668").Value = 969
$ws.Range("I671").Value = 884
$ws.Range("I676").Value = 681
$ws.Range("I687").Value = 828
$ws.Range("I707").Value = 562
$ws.Range("I712").Value = 883
$ws.Range("I733").Value = 493
$ws.Range("I736").Value = 919
$ws.Range("I746").Value = 561
$ws.Range("I757").Value = 511
$ws.Range("I761").Value = 641
$ws.Range("I764").Value = 953
$ws.Range("I773").Value = 512
$ws.Range("I775").Value = 901
$ws.Range("I776").Value = 865
$ws.Range("I785").Value = 680
$ws.Range("I787").Value = 867
$ws.Range("I797").Value = 952
$ws.Range("I798").Value = 633
$ws.Range("I813").Value = 950
$ws.Range("I823").Value = 932
$ws.Range("I827").Value = 614
$ws.Range("I828").Value = 560
$ws.Range("I841").Value = 743
$ws.Range("I849").Value = 797
$ws.Range("I856").Value = 742
$ws.Range("I858").Value = 778
$ws.Range("I867").Value = 791
$ws.Range("I868").Value = 938
$ws.Range("I876").Value = 988
$ws.Range("I882").Value = 317
$ws.Range("I888").Value = 937
$ws.Range("I894").Value = 994
$ws.Range("I898").Value = 494
$ws.Range("I900").Value = 913
$ws.Range("I904").Value = 749
$ws.Range("I909").Value = 818
$ws.Range("I918").Value = 908
$ws.Range("I928").Value = 728
$ws.Range("I934").Value = 957
$ws.Range("I937").Value = 811
$ws.Range("I943").Value = 873
$ws.Range("I944").Value = 593
$ws.Range("I947").Value = 379
$ws.Range("I957").Value = 634
$ws.Range("I963").Value = 912
$ws.Range("I964").Value = 961
$ws.Range("I965").Value = 983
$ws.Range("I967").Value = 1031
